# C5-PowerPoint.pptx edit
#
# 1) The table on slide 6 (graphicFrame "Google Shape;127;p18") switches
#    from the deck's custom "Table_0" style to the built-in PowerPoint
#    table style {CFE34F9A-54EE-4E0A-B570-CC87A9FF570A}.
#
# 2) The presentation's design theme is changed from "Integral" to the
#    standard Office theme (the 12 theme colours used by the slide master
#    change from the Integral palette to the default Office palette).

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 6 -------------------------------------
$slide6 = $p.Slides.Item(6)
for ($i = 1; $i -le $slide6.Shapes.Count; $i++) {
    $shp = $slide6.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{CFE34F9A-54EE-4E0A-B570-CC87A9FF570A}")
    }
}

# --- 2. Switch the deck's theme colours from Integral to Office -----------
$tcs = $p.Slides.Item(1).ThemeColorScheme
$tcs.Item(1).RGB  = 0         # Dark 1    - 000000
$tcs.Item(2).RGB  = 16777215  # Light 1   - FFFFFF
$tcs.Item(3).RGB  = 6968388   # Dark 2    - 44546A
$tcs.Item(4).RGB  = 15132391  # Light 2   - E7E6E6
$tcs.Item(5).RGB  = 13998939  # Accent 1  - 5B9BD5
$tcs.Item(6).RGB  = 3243501   # Accent 2  - ED7D31
$tcs.Item(7).RGB  = 10855845  # Accent 3  - A5A5A5
$tcs.Item(8).RGB  = 49407     # Accent 4  - FFC000
$tcs.Item(9).RGB  = 12874308  # Accent 5  - 4472C4
$tcs.Item(10).RGB = 4697456   # Accent 6  - 70AD47
$tcs.Item(11).RGB = 12673797  # Hyperlink - 0563C1
$tcs.Item(12).RGB = 7491477   # Followed Hyperlink - 954F72
